$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("D4").Value = 44405
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 1200
$ws.Range("O4").Value = 1200
$ws.Range("P4").Value = 1200
$ws.Range("S4").Value = 1200

# Row 6
$ws.Range("D6").Value = 44418
$ws.Range("M6").Value = 40

# Row 7
$ws.Range("D7").Value = 44343
$ws.Range("N7").Value = 1300
$ws.Range("O7").Value = 1300
$ws.Range("P7").Value = 1300
$ws.Range("S7").Value = 1300

# Row 9
$ws.Range("D9").Value = 44432
$ws.Range("M9").Value = 30

# Row 10
$ws.Range("D10").Value = 44438
$ws.Range("M10").Value = 60

# Row 11
$ws.Range("D11").Value = 44435
$ws.Range("M11").Value = 130
